# edit.ps1 — applies the "add 2022-Q3 data" change:
#  1. Inserts a new worksheet "2022-Q3" right after "总计" (before "2022-Q2"),
#     populated with the fund-holdings table for that quarter.
#  2. Prepends a corresponding summary row to the "总计" sheet, shifting the
#     existing quarters down by one row (indices/row numbers +1).

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------------
# 1) "总计" sheet: shift existing rows down by re-writing values (row 2 becomes
#    the new 2022-Q3 entry, rows 3-8 hold what used to be rows 2-7).
# ---------------------------------------------------------------------------

$totalRows = @(
  @(0, "2022-Q3", 31, 5.78),
  @(1, "2022-Q2", 29, 3.81),
  @(2, "2022-Q1", 6, 0.63),
  @(3, "2021-Q4", 20, 3.88),
  @(4, "2021-Q3", 29, 6.34),
  @(5, "2021-Q1", 4, 0.32),
  @(6, "2020-Q4", 3, 0.44)
)

for ($i = 0; $i -lt $totalRows.Length; $i++) {
  $r = 2 + $i
  $entry = $totalRows[$i]
  $totalSheet.Cells.Item($r, 1).Value = $entry[0]
  $totalSheet.Cells.Item($r, 2).Value = $entry[1]
  $totalSheet.Cells.Item($r, 3).Value = $entry[2]
  $totalSheet.Cells.Item($r, 4).Value = $entry[3]
}

# Row 8 (2020-Q4) is a brand-new row — give its A cell the same style ("s=2",
# bold/centered/bordered) the other index cells in column A already use, by
# copying the format from the row directly above it.
$totalSheet.Cells.Item(7, 1).Copy()
$totalSheet.Cells.Item(8, 1).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) New "2022-Q3" worksheet, inserted right before "2022-Q2".
# ---------------------------------------------------------------------------

$ws = $wb.Worksheets.Add($q2Sheet)
$ws.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, 2 + $i).Value = $headers[$i]
}
# Header style: bold/centered/bordered — same "s=2" style used on row 1 of
# every other quarter sheet. Copy it over from the "总计" sheet header cell.
$totalSheet.Cells.Item(1, 2).Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

# Tab-separated data block, one row per fund. Each field is encoded as
# "<kind>|<value>" where kind is S (text) or N (number).
$dataBlock = @"
S|001224	S|中邮新思路灵活配置混合	S|27.22	S|74.71	S|4.12	S|1.1215	N|6
S|001245	S|工银生态环境股票A	S|37.36	S|88.46	S|2.21	S|0.8257	N|10
S|007777	S|中邮研究精选混合	S|14.89	S|70.15	S|4.67	S|0.6954	N|3
S|012975	S|西部利得碳中和混合A	S|10.84	S|92.98	S|4.78	S|0.5182	N|3
S|002620	S|中邮未来新蓝筹灵活配置混合	S|11.81	S|83.48	S|3.93	S|0.4641	N|7
S|001910	S|泰康新机遇灵活配置混合	S|15.48	S|82.44	S|2.55	S|0.3947	N|10
S|007040	S|新疆前海联合泳隆灵活配置混合C	S|6.78	S|91.50	S|3.90	S|0.2644	N|7
S|012976	S|西部利得碳中和混合C	S|5.31	S|92.98	S|4.78	S|0.2538	N|3
S|008980	S|中邮科技创新精选混合A	S|4.75	S|88.40	S|4.70	S|0.2232	N|6
S|009490	S|泰康科技创新一年定期开放混合	S|2.44	S|79.62	S|7.25	S|0.1769	N|1
S|673060	S|西部利得景瑞灵活配置混合A	S|3.97	S|93.10	S|4.13	S|0.1640	N|7
S|011001	S|中邮兴荣价值一年持有期混合	S|5.15	S|40.76	S|3.00	S|0.1545	N|4
S|008981	S|中邮科技创新精选混合C	S|2.66	S|88.40	S|4.70	S|0.1250	N|6
S|011793	S|建信智能汽车股票	S|4.86	S|91.44	S|2.29	S|0.1113	N|10
S|009258	S|西部利得景瑞灵活配置混合C	S|1.32	S|93.10	S|4.13	S|0.0545	N|7
S|005933	S|新疆前海联合先进制造灵活配置混合A	S|0.95	S|92.14	S|4.42	S|0.0420	N|9
S|002935	S|泰康恒泰回报灵活配置混合C	S|2.53	S|22.01	S|1.64	S|0.0415	N|5
S|004128	S|新疆前海联合泳隆灵活配置混合A	S|1.02	S|91.50	S|3.90	S|0.0398	N|7
S|014938	S|同泰产业升级混合A	S|1.01	S|61.58	S|2.68	S|0.0271	N|9
S|002934	S|泰康恒泰回报灵活配置混合A	S|1.12	S|22.01	S|1.64	S|0.0184	N|5
S|166109	S|信澳量化先锋混合（LOF）A	S|0.79	S|88.99	S|2.32	S|0.0183	N|10
S|007770	S|同泰开泰混合A	S|0.40	S|85.52	S|2.75	S|0.0110	N|7
S|519961	S|长信利广灵活配置混合A	S|0.74	S|26.85	S|1.29	S|0.0095	N|7
S|015002	S|工银生态环境股票C	S|0.26	S|88.46	S|2.21	S|0.0057	N|10
S|501002	S|长信价值优选混合	S|0.39	S|93.83	S|1.44	S|0.0056	N|6
S|007771	S|同泰开泰混合C	S|0.17	S|85.52	S|2.75	S|0.0047	N|7
S|005934	S|新疆前海联合先进制造灵活配置混合C	S|0.10	S|92.14	S|4.42	S|0.0044	N|9
S|008890	S|中邮价值优选一年定期开放灵活配置混合	S|0.12	S|61.34	S|3.14	S|0.0038	N|8
S|166110	S|信澳量化先锋混合（LOF）C	S|0.11	S|88.99	S|2.32	S|0.0026	N|10
S|519960	S|长信利广灵活配置混合C	S|0.03	S|26.85	S|1.29	S|0.0004	N|7
S|014939	S|同泰产业升级混合C	S|0.00	S|61.58	S|2.68	N|0	N|9
"@

$lines = $dataBlock -split "`n"
$rowIdx = 2
foreach ($line in $lines) {
  $trimmed = $line.Trim()
  if ($trimmed.Length -eq 0) { continue }
  $fields = $trimmed -split "`t"

  # Column A: 0-based row index, numeric, same "s=2" style as the other sheets.
  $ws.Cells.Item($rowIdx, 1).Value = $rowIdx - 2

  for ($c = 0; $c -lt $fields.Length; $c++) {
    $kv = $fields[$c] -split "\|", 2
    $kind = $kv[0]
    $val = $kv[1]
    $cell = $ws.Cells.Item($rowIdx, 2 + $c)
    if ($kind -eq "N") {
      $cell.Value = [double]$val
    } else {
      # Force text storage (matches the source data, which keeps these as
      # strings even though they look numeric) without leaving a stray
      # "Text" number-format style behind.
      $cell.NumberFormat = "@"
      $cell.Value = $val
      $cell.ClearFormats()
    }
  }
  $rowIdx++
}

# Column A style (index column), matches the "s=2" style used elsewhere.
$totalSheet.Cells.Item(2, 1).Copy()
$ws.Range("A2:A32").PasteSpecial(-4122)
